$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 5-7 (task 1.3, 1.4, 1.5) as "done" in the Status column (E),
# matching the already-"done" rows 3 and 4.
$ws.Range("E5").Value = "done"
$ws.Range("E6").Value = "done"
$ws.Range("E7").Value = "done"

# Move the active selection to D14, as last left by the editor.
$ws.Range("D14").Select()
